# Scheduled market-data refresh: update currentAveragePrice* / Leve profit
# columns (H-N) across the profession sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 461415.4
$ws.Range("J17").Value = 461415.4
$ws.Range("L17").Value = 1384246.2
$ws.Range("N17").Value = -1384582.2

$ws.Range("H121").Value = 1779.8
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1779.8
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 5339.4
$ws.Range("N121").Value = -8833.4
$ws.Range("M121").ClearContents()

$ws.Range("H137").Value = 3227528.5
$ws.Range("I137").Value = 6251071.5
$ws.Range("K137").Value = 18753214.5
$ws.Range("M137").Value = -18750664.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 133467690
$ws.Range("I61").Value = 83417780
$ws.Range("K61").Value = 83417780
$ws.Range("M61").Value = -83417568

$ws.Range("H74").Value = 11455558
$ws.Range("I74").Value = 14765709
$ws.Range("J74").Value = 201046.6
$ws.Range("K74").Value = 14765709
$ws.Range("L74").Value = 201046.6
$ws.Range("M74").Value = -14764835
$ws.Range("N74").Value = -202794.6

$ws.Range("H77").Value = 11455558
$ws.Range("I77").Value = 14765709
$ws.Range("J77").Value = 201046.6
$ws.Range("K77").Value = 73828545
$ws.Range("L77").Value = 1005233
$ws.Range("M77").Value = -73824177
$ws.Range("N77").Value = -1013969

$ws.Range("H136").Value = 133467690
$ws.Range("I136").Value = 83417780
$ws.Range("K136").Value = 250253340
$ws.Range("M136").Value = -250250790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 789
$ws.Range("J64").Value = 838.3333
$ws.Range("L64").Value = 838.3333
$ws.Range("N64").Value = -1288.3333

$ws.Range("H67").Value = 789
$ws.Range("J67").Value = 838.3333
$ws.Range("L67").Value = 838.3333
$ws.Range("N67").Value = -2398.3333

$ws.Range("H134").Value = 2809.85
$ws.Range("I134").Value = 3016.0625
$ws.Range("K134").Value = 9048.1875
$ws.Range("M134").Value = -6513.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2625.3062
$ws.Range("I31").Value = 1042.4166
$ws.Range("K31").Value = 1042.4166
$ws.Range("M31").Value = -747.4166

$ws.Range("H34").Value = 2625.3062
$ws.Range("I34").Value = 1042.4166
$ws.Range("K34").Value = 1042.4166
$ws.Range("M34").Value = -840.4166

$ws.Range("H58").Value = 48782210
$ws.Range("I58").Value = 41667788
$ws.Range("J58").Value = 58826100
$ws.Range("K58").Value = 41667788
$ws.Range("L58").Value = 58826100
$ws.Range("M58").Value = -41667585
$ws.Range("N58").Value = -58826506

$ws.Range("H86").Value = 4041.5908
$ws.Range("I86").Value = 3977.1177
$ws.Range("J86").Value = 4260.8
$ws.Range("K86").Value = 3977.1177
$ws.Range("L86").Value = 4260.8
$ws.Range("M86").Value = -2854.1177
$ws.Range("N86").Value = -6506.8

$ws.Range("H89").Value = 4041.5908
$ws.Range("I89").Value = 3977.1177
$ws.Range("J89").Value = 4260.8
$ws.Range("K89").Value = 19885.5885
$ws.Range("L89").Value = 21304
$ws.Range("M89").Value = -14269.5885
$ws.Range("N89").Value = -32536

$ws.Range("H100").Value = 38000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 38000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 38000
$ws.Range("N100").Value = -40164
$ws.Range("M100").ClearContents()

$ws.Range("H134").Value = 43496.5
$ws.Range("I134").Value = 2005.619
$ws.Range("J134").Value = 217758.2
$ws.Range("K134").Value = 6016.857
$ws.Range("L134").Value = 653274.6000000001
$ws.Range("M134").Value = -3481.857
$ws.Range("N134").Value = -658344.6000000001

$ws.Range("H136").Value = 48782210
$ws.Range("I136").Value = 41667788
$ws.Range("J136").Value = 58826100
$ws.Range("K136").Value = 125003364
$ws.Range("L136").Value = 176478300
$ws.Range("M136").Value = -125000814
$ws.Range("N136").Value = -176483400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1382.6471
$ws.Range("J132").Value = 1710.5
$ws.Range("L132").Value = 15394.5
$ws.Range("N132").Value = -20454.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1828.3
$ws.Range("I126").Value = 1130.2222
$ws.Range("J126").Value = 2399.4546
$ws.Range("K126").Value = 3390.6666
$ws.Range("L126").Value = 7198.3638
$ws.Range("M126").Value = -920.6665999999996
$ws.Range("N126").Value = -12138.3638

$ws.Range("H133").Value = 49250
$ws.Range("J133").Value = 49250
$ws.Range("L133").Value = 49250
$ws.Range("N133").Value = -59370

$ws.Range("H135").Value = 29975.562
$ws.Range("I135").Value = 30709
$ws.Range("J135").Value = 29926.666
$ws.Range("K135").Value = 30709
$ws.Range("L135").Value = 29926.666
$ws.Range("M135").Value = -25639
$ws.Range("N135").Value = -40066.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 93361.73
$ws.Range("I132").Value = 668
$ws.Range("J132").Value = 128121.875
$ws.Range("K132").Value = 2004
$ws.Range("L132").Value = 384365.625
$ws.Range("M132").Value = 526
$ws.Range("N132").Value = -389425.625

$ws.Range("H136").Value = 223566.56
$ws.Range("I136").Value = 334433.34
$ws.Range("J136").Value = 168133.17
$ws.Range("K136").Value = 1003300.02
$ws.Range("L136").Value = 504399.51
$ws.Range("M136").Value = -1000750.02
$ws.Range("N136").Value = -509499.51

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4239.143
$ws.Range("I62").Value = 3875
$ws.Range("J62").Value = 4724.6665
$ws.Range("K62").Value = 3875
$ws.Range("L62").Value = 4724.6665
$ws.Range("M62").Value = -3251
$ws.Range("N62").Value = -5972.6665

$ws.Range("H65").Value = 4239.143
$ws.Range("I65").Value = 3875
$ws.Range("J65").Value = 4724.6665
$ws.Range("K65").Value = 19375
$ws.Range("L65").Value = 23623.3325
$ws.Range("M65").Value = -16255
$ws.Range("N65").Value = -29863.3325

$ws.Range("H132").Value = 107113.79
$ws.Range("I132").Value = 144114.86
$ws.Range("J132").Value = 85529.836
$ws.Range("K132").Value = 432344.58
$ws.Range("L132").Value = 256589.508
$ws.Range("M132").Value = -429814.58
$ws.Range("N132").Value = -261649.508

$ws.Range("H136").Value = 43218.668
$ws.Range("I136").Value = 25906.45
$ws.Range("J136").Value = 129779.75
$ws.Range("K136").Value = 77719.35000000001
$ws.Range("L136").Value = 389339.25
$ws.Range("M136").Value = -75169.35000000001
$ws.Range("N136").Value = -394439.25
